# Add a new "Upcoming Semantic Interop Events" slide right before the final
# "Any Other Business" slide (i.e. insert at position 23, pushing the old
# slide 23 to become slide 24).

$p = $ppt.ActivePresentation

# ppLayoutText (2) => "Title and Content" CustomLayout on this master.
$newSlide = $p.Slides.Add(23, 2)

# --- Title -----------------------------------------------------------
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Upcoming Semantic Interop Events"

# --- Body / content placeholder ---------------------------------------
$bodyShape = $newSlide.Shapes.Item(2)
$body = $bodyShape.TextFrame.TextRange

$body.Text = "W3C Web of Things Plugfest"
$body.InsertAfter("`rSouth Korea, June 30th and July 1st") | Out-Null
$body.InsertAfter("`rWISHI Plugfest/Hackathon at IETF 102") | Out-Null
$body.InsertAfter("`rMontreal, July 14th and 15th") | Out-Null

# Superscript the ordinal suffixes "th" (June 30th) and "st" (July 1st).
$th = $body.Characters(48, 2)
$th.Font.Superscript = $true

$st = $body.Characters(61, 2)
$st.Font.Superscript = $true

# Second-level bullets for the two "location/date" lines.
$body.Paragraphs(2, 1).IndentLevel = 2
$body.Paragraphs(4, 1).IndentLevel = 2

# Reposition/resize the content placeholder to match the authored layout.
$bodyShape.Left = 49.5
$bodyShape.Top = 160.93185
$bodyShape.Width = 621.0
$bodyShape.Height = 342.62508
